$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 727.9231
$ws.Cells.Item(19, 9).Value = 696.6667
$ws.Cells.Item(19, 10).Value = 737.3
$ws.Cells.Item(19, 11).Value = 696.6667
$ws.Cells.Item(19, 12).Value = 737.3
$ws.Cells.Item(19, 13).Value = -521.6667
$ws.Cells.Item(19, 14).Value = -1087.3
$ws.Cells.Item(51, 8).Value = 3770
$ws.Cells.Item(51, 9).Value = 3000
$ws.Cells.Item(51, 10).Value = 4100
$ws.Cells.Item(51, 11).Value = 3000
$ws.Cells.Item(51, 12).Value = 4100
$ws.Cells.Item(51, 13).Value = -2516
$ws.Cells.Item(51, 14).Value = -5068
$ws.Cells.Item(121, 8).Value = 1529.1666
$ws.Cells.Item(121, 10).Value = 2293.3333
$ws.Cells.Item(121, 12).Value = 6879.999899999999
$ws.Cells.Item(121, 14).Value = -10373.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10744.72
$ws.Cells.Item(32, 9).Value = 10885.046
$ws.Cells.Item(32, 10).Value = 9715.666999999999
$ws.Cells.Item(32, 11).Value = 10885.046
$ws.Cells.Item(32, 12).Value = 9715.666999999999
$ws.Cells.Item(32, 13).Value = -10598.046
$ws.Cells.Item(32, 14).Value = -10289.667
$ws.Cells.Item(42, 8).Value = 50000
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 1593.4324
$ws.Cells.Item(61, 9).Value = 1381.2759
$ws.Cells.Item(61, 10).Value = 2362.5
$ws.Cells.Item(61, 11).Value = 1381.2759
$ws.Cells.Item(61, 12).Value = 2362.5
$ws.Cells.Item(61, 13).Value = -1169.2759
$ws.Cells.Item(61, 14).Value = -2786.5
$ws.Cells.Item(74, 8).Value = 1641.4286
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 1641.4286
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(117, 8).Value = 24500
$ws.Cells.Item(117, 10).Value = 24500
$ws.Cells.Item(117, 12).Value = 24500
$ws.Cells.Item(117, 14).Value = -33678
$ws.Cells.Item(137, 8).Value = 1352.1875
$ws.Cells.Item(137, 9).Value = 1309
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 3927
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 13).Value = -1377
$ws.Cells.Item(137, 14).Value = -11100
$ws.Cells.Item(141, 8).Value = 4643.1143
$ws.Cells.Item(141, 9).Value = 1621.8667
$ws.Cells.Item(141, 10).Value = 22770.6
$ws.Cells.Item(141, 11).Value = 4865.6001
$ws.Cells.Item(141, 12).Value = 68311.79999999999
$ws.Cells.Item(141, 13).Value = 314.3999000000003
$ws.Cells.Item(141, 14).Value = -78671.79999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2083.5454
$ws.Cells.Item(99, 9).Value = 1501.6666
$ws.Cells.Item(99, 10).Value = 2781.8
$ws.Cells.Item(99, 11).Value = 1501.6666
$ws.Cells.Item(99, 12).Value = 2781.8
$ws.Cells.Item(99, 13).Value = -3.666600000000017
$ws.Cells.Item(99, 14).Value = -5777.8
$ws.Cells.Item(132, 8).Value = 4140.959
$ws.Cells.Item(132, 9).Value = 5193.1377
$ws.Cells.Item(132, 10).Value = 2615.3
$ws.Cells.Item(132, 11).Value = 15579.4131
$ws.Cells.Item(132, 12).Value = 7845.900000000001
$ws.Cells.Item(132, 13).Value = -13049.4131
$ws.Cells.Item(132, 14).Value = -12905.9
$ws.Cells.Item(136, 8).Value = 1593.4324
$ws.Cells.Item(136, 9).Value = 1381.2759
$ws.Cells.Item(136, 10).Value = 2362.5
$ws.Cells.Item(136, 11).Value = 4143.8277
$ws.Cells.Item(136, 12).Value = 7087.5
$ws.Cells.Item(136, 13).Value = -1593.8277
$ws.Cells.Item(136, 14).Value = -12187.5
$ws.Cells.Item(140, 8).Value = 54405
$ws.Cells.Item(140, 10).Value = 54405
$ws.Cells.Item(140, 12).Value = 54405
$ws.Cells.Item(140, 14).Value = -64765

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 4
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 135
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 48199.6
$ws.Cells.Item(20, 10).Value = 48199.6
$ws.Cells.Item(20, 12).Value = 48199.6
$ws.Cells.Item(20, 14).Value = -48671.6
$ws.Cells.Item(30, 8).Value = 48199.6
$ws.Cells.Item(30, 10).Value = 48199.6
$ws.Cells.Item(30, 12).Value = 48199.6
$ws.Cells.Item(30, 14).Value = -48381.6
$ws.Cells.Item(31, 8).Value = 1767.3043
$ws.Cells.Item(31, 9).Value = 1696.8572
$ws.Cells.Item(31, 10).Value = 2507
$ws.Cells.Item(31, 11).Value = 1696.8572
$ws.Cells.Item(31, 12).Value = 2507
$ws.Cells.Item(31, 13).Value = -1401.8572
$ws.Cells.Item(31, 14).Value = -3097
$ws.Cells.Item(34, 8).Value = 1767.3043
$ws.Cells.Item(34, 9).Value = 1696.8572
$ws.Cells.Item(34, 10).Value = 2507
$ws.Cells.Item(34, 11).Value = 1696.8572
$ws.Cells.Item(34, 12).Value = 2507
$ws.Cells.Item(34, 13).Value = -1494.8572
$ws.Cells.Item(34, 14).Value = -2911
$ws.Cells.Item(58, 8).Value = 700456.3
$ws.Cells.Item(58, 9).Value = 1030243.9
$ws.Cells.Item(58, 10).Value = 2082.5293
$ws.Cells.Item(58, 11).Value = 1030243.9
$ws.Cells.Item(58, 12).Value = 2082.5293
$ws.Cells.Item(58, 13).Value = -1030040.9
$ws.Cells.Item(58, 14).Value = -2488.5293
$ws.Cells.Item(68, 8).Value = 32000
$ws.Cells.Item(68, 10).Value = 32000
$ws.Cells.Item(68, 12).Value = 32000
$ws.Cells.Item(68, 14).Value = -33498
$ws.Cells.Item(71, 8).Value = 32000
$ws.Cells.Item(71, 10).Value = 32000
$ws.Cells.Item(71, 12).Value = 96000
$ws.Cells.Item(71, 14).Value = -103488
$ws.Cells.Item(74, 8).Value = 33285.715
$ws.Cells.Item(74, 10).Value = 33285.715
$ws.Cells.Item(74, 12).Value = 33285.715
$ws.Cells.Item(74, 14).Value = -35033.715
$ws.Cells.Item(77, 8).Value = 33285.715
$ws.Cells.Item(77, 10).Value = 33285.715
$ws.Cells.Item(77, 12).Value = 99857.14499999999
$ws.Cells.Item(77, 14).Value = -108593.145
$ws.Cells.Item(107, 8).Value = 417.3684
$ws.Cells.Item(107, 9).Value = 351.91666
$ws.Cells.Item(107, 10).Value = 529.5714
$ws.Cells.Item(107, 11).Value = 351.91666
$ws.Cells.Item(107, 12).Value = 529.5714
$ws.Cells.Item(107, 13).Value = 1568.08334
$ws.Cells.Item(107, 14).Value = -4369.5714
$ws.Cells.Item(128, 8).Value = 48199.6
$ws.Cells.Item(128, 10).Value = 48199.6
$ws.Cells.Item(128, 12).Value = 48199.6
$ws.Cells.Item(128, 14).Value = -58159.6
$ws.Cells.Item(134, 8).Value = 1897.1091
$ws.Cells.Item(134, 9).Value = 1557.4783
$ws.Cells.Item(134, 10).Value = 3633
$ws.Cells.Item(134, 11).Value = 4672.4349
$ws.Cells.Item(134, 12).Value = 10899
$ws.Cells.Item(134, 13).Value = -2137.4349
$ws.Cells.Item(134, 14).Value = -15969

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 700456.3
$ws.Cells.Item(136, 9).Value = 1030243.9
$ws.Cells.Item(136, 10).Value = 2082.5293
$ws.Cells.Item(136, 11).Value = 3090731.7
$ws.Cells.Item(136, 12).Value = 6247.5879
$ws.Cells.Item(136, 13).Value = -3088181.7
$ws.Cells.Item(136, 14).Value = -11347.5879

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 20080120
$ws.Cells.Item(7, 9).Value = 33333534
$ws.Cells.Item(7, 10).Value = 200000
$ws.Cells.Item(7, 11).Value = 33333534
$ws.Cells.Item(7, 12).Value = 200000
$ws.Cells.Item(7, 13).Value = -33333422
$ws.Cells.Item(7, 14).Value = -200224
$ws.Cells.Item(8, 8).Value = 20080120
$ws.Cells.Item(8, 9).Value = 33333534
$ws.Cells.Item(8, 10).Value = 200000
$ws.Cells.Item(8, 11).Value = 33333534
$ws.Cells.Item(8, 12).Value = 200000
$ws.Cells.Item(8, 13).Value = -33333395
$ws.Cells.Item(8, 14).Value = -200278
$ws.Cells.Item(109, 8).Value = 9648.947
$ws.Cells.Item(109, 10).Value = 9648.947
$ws.Cells.Item(109, 12).Value = 9648.947
$ws.Cells.Item(109, 14).Value = -11728.947
$ws.Cells.Item(126, 8).Value = 4218.1816
$ws.Cells.Item(126, 9).Value = 4062.5
$ws.Cells.Item(126, 10).Value = 4633.3335
$ws.Cells.Item(126, 11).Value = 12187.5
$ws.Cells.Item(126, 12).Value = 13900.0005
$ws.Cells.Item(126, 13).Value = -9717.5
$ws.Cells.Item(126, 14).Value = -18840.0005
$ws.Cells.Item(133, 8).Value = 3655.818
$ws.Cells.Item(133, 9).Value = 1630
$ws.Cells.Item(133, 11).Value = 4890
$ws.Cells.Item(133, 13).Value = 170
$ws.Cells.Item(138, 8).Value = 2584.4
$ws.Cells.Item(138, 9).Value = 1095
$ws.Cells.Item(138, 10).Value = 4286.5713
$ws.Cells.Item(138, 11).Value = 3285
$ws.Cells.Item(138, 12).Value = 12859.7139
$ws.Cells.Item(138, 13).Value = 1855
$ws.Cells.Item(138, 14).Value = -23139.7139
$ws.Cells.Item(140, 8).Value = 3813.4412
$ws.Cells.Item(140, 9).Value = 964.7368
$ws.Cells.Item(140, 10).Value = 7421.8
$ws.Cells.Item(140, 11).Value = 2894.2104
$ws.Cells.Item(140, 12).Value = 22265.4
$ws.Cells.Item(140, 13).Value = 2285.7896
$ws.Cells.Item(140, 14).Value = -32625.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1916.659
$ws.Cells.Item(132, 9).Value = 1249.7742
$ws.Cells.Item(132, 10).Value = 3506.923
$ws.Cells.Item(132, 11).Value = 3749.3226
$ws.Cells.Item(132, 12).Value = 10520.769
$ws.Cells.Item(132, 13).Value = -1219.3226
$ws.Cells.Item(132, 14).Value = -15580.769
$ws.Cells.Item(136, 8).Value = 36863
$ws.Cells.Item(136, 10).Value = 36863
$ws.Cells.Item(136, 12).Value = 110589
$ws.Cells.Item(136, 14).Value = -115689

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 1000000
$ws.Cells.Item(5, 9).Value = 1000000
$ws.Cells.Item(5, 11).Value = 1000000
$ws.Cells.Item(5, 13).Value = -999888
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 5907.143
$ws.Cells.Item(132, 9).Value = 5972.909
$ws.Cells.Item(132, 11).Value = 17918.727
$ws.Cells.Item(132, 13).Value = -15388.727
$ws.Cells.Item(136, 8).Value = 24635990
$ws.Cells.Item(136, 9).Value = 28572760
$ws.Cells.Item(136, 10).Value = 1671500.9
$ws.Cells.Item(136, 11).Value = 85718280
$ws.Cells.Item(136, 12).Value = 5014502.699999999
$ws.Cells.Item(136, 13).Value = -85715730
$ws.Cells.Item(136, 14).Value = -5019602.699999999
